$wb = $excel.ActiveWorkbook

# --- Sheet "DBD" (sheet1): UspErrorLog field/table layout ---
$ws1 = $wb.Worksheets.Item("DBD")

# Index2 definition (row 6) now references the new JobTxSeq column
$ws1.Cells.Item(6, 3).Value2 = "JobTxSeq"

# New field row (SEQ 13): JobTxSeq column added to the UspErrorLog table
$ws1.Cells.Item(21, 1).Value2 = 13
$ws1.Cells.Item(21, 2).Value2 = "JobTxSeq"
$ws1.Cells.Item(21, 3).Value2 = "啟動批次的交易序號"
$ws1.Cells.Item(21, 4).Value2 = "VARCHAR2"
$ws1.Cells.Item(21, 5).Value2 = 20

# --- Sheet "DBS" (sheet2): query function definitions ---
$ws2 = $wb.Worksheets.Item("DBS")

# New lookup function: findByJobTxSeq
$ws2.Cells.Item(4, 1).Value2 = "findByJobTxSeq"
$ws2.Cells.Item(4, 2).Value2 = "JobTxSeq = "
$ws2.Cells.Item(4, 3).Value2 = "LogDate DESC,LogTime DESC,UspName DESC"

$ws2.Range("A5").Select()

# Re-activate the DBD tab and leave the selection where the edit happened
# (matches the author's saved view: DBD tab selected, C6 highlighted)
$ws1.Activate()
$ws1.Range("C6").Select()

$wb.Save()
